# feat: add GPA record access
#
# Mark the "Configure database interface for student GPA" row (and the two
# related rows around it) as done by using the same "2 down" marker that is
# already used elsewhere in column D of this checklist, instead of a plain
# numeric "2".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D21").Value = "2 down"
$ws.Range("D22").Value = "2 down"
$ws.Range("D30").Value = "2 down"

# Move the active selection to G23, matching where the edit was made.
$ws.Range("G23").Select()
